# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect newly scraped counts (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---------------------------------------------------------
$wsExhibit = $wb.Worksheets.Item("展览")

$exhibitChanges = @{
    2  = 1111
    3  = 434
    4  = 1524
    5  = 8810
    8  = 656
    9  = 302
    12 = 21
    13 = 3670
    16 = 87
    17 = 2817
    18 = 0
    19 = 1126
    20 = 317
    21 = 216
    22 = 2460
    23 = 80
}

foreach ($row in $exhibitChanges.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitChanges[$row]
}

# --- Sheet "全部类型" ------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")

$allChanges = @{
    2  = 1111
    3  = 434
    4  = 1524
    5  = 8810
    8  = 656
    9  = 302
    12 = 21
    13 = 3670
    16 = 87
    17 = 2818
    18 = 151
    19 = 1126
    20 = 317
    21 = 216
    22 = 2460
    24 = 80
}

foreach ($row in $allChanges.Keys) {
    $wsAll.Range("F$row").Value = $allChanges[$row]
}
